# Update res_bus vm_pu values for Case_2_25 (380 kV case)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 2).Value = 1.02
$ws.Cells.Item(2, 3).Value = 1.028417834686724
$ws.Cells.Item(2, 4).Value = 1.031885696139161
$ws.Cells.Item(2, 5).Value = 0.9926147277508489
$ws.Cells.Item(2, 6).Value = 1.026968078667254
$ws.Cells.Item(2, 9).Value = 1.03358128967705
$ws.Cells.Item(2, 10).Value = 1.033570265693573
$ws.Cells.Item(2, 11).Value = 1.034692521117107
$ws.Cells.Item(2, 12).Value = 0.9955398523336033
$ws.Cells.Item(2, 13).Value = 1.029789178236294
$ws.Cells.Item(2, 14).Value = 1.015178337871

$ws.Cells.Item(3, 2).Value = 1.02
$ws.Cells.Item(3, 3).Value = 1.02939674924391
$ws.Cells.Item(3, 4).Value = 1.032597889562139
$ws.Cells.Item(3, 5).Value = 0.9936372048519304
$ws.Cells.Item(3, 6).Value = 1.028580580264547
$ws.Cells.Item(3, 9).Value = 1.033809814566775
$ws.Cells.Item(3, 10).Value = 1.034189629816018
$ws.Cells.Item(3, 11).Value = 1.035213630079863
$ws.Cells.Item(3, 12).Value = 0.9963617723202692
$ws.Cells.Item(3, 13).Value = 1.031207129397451
$ws.Cells.Item(3, 14).Value = 1.015385868337456

$ws.Cells.Item(4, 2).Value = 1.02
$ws.Cells.Item(4, 3).Value = 1.030029705008401
$ws.Cells.Item(4, 4).Value = 1.033058353155536
$ws.Cells.Item(4, 5).Value = 0.9942998659930995
$ws.Cells.Item(4, 6).Value = 1.029623448226028
$ws.Cells.Item(4, 9).Value = 1.03395630332032
$ws.Cells.Item(4, 10).Value = 1.034589349338294
$ws.Cells.Item(4, 11).Value = 1.035549779809773
$ws.Cells.Item(4, 12).Value = 0.9968940712668345
$ws.Cells.Item(4, 13).Value = 1.032123646250682
$ws.Cells.Item(4, 14).Value = 1.015519740192953

$ws.Cells.Item(5, 2).Value = 1.02
$ws.Cells.Item(5, 3).Value = 1.03029568876282
$ws.Cells.Item(5, 4).Value = 1.033251842822806
$ws.Cells.Item(5, 5).Value = 0.9945786998346017
$ws.Cells.Item(5, 6).Value = 1.030061750039104
$ws.Cells.Item(5, 9).Value = 1.034017556545188
$ws.Cells.Item(5, 10).Value = 1.034757140690339
$ws.Cells.Item(5, 11).Value = 1.035690847972759
$ws.Cells.Item(5, 12).Value = 0.997117960005301
$ws.Cells.Item(5, 13).Value = 1.032508718104115
$ws.Cells.Item(5, 14).Value = 1.015575920937035

$ws.Cells.Item(6, 2).Value = 1.02
$ws.Cells.Item(6, 3).Value = 1.030340342135719
$ws.Cells.Item(6, 4).Value = 1.033284325360049
$ws.Cells.Item(6, 5).Value = 0.9946255319796338
$ws.Cells.Item(6, 6).Value = 1.030135336004473
$ws.Cells.Item(6, 9).Value = 1.034027821848717
$ws.Cells.Item(6, 10).Value = 1.034785298915146
$ws.Cells.Item(6, 11).Value = 1.035714519337901
$ws.Cells.Item(6, 12).Value = 0.9971555583673453
$ws.Cells.Item(6, 13).Value = 1.032573359934501
$ws.Cells.Item(6, 14).Value = 1.015585348132004

$ws.Cells.Item(7, 2).Value = 1.02
$ws.Cells.Item(7, 3).Value = 1.030033259530606
$ws.Cells.Item(7, 4).Value = 1.033060938924425
$ws.Cells.Item(7, 5).Value = 0.9943035907982488
$ws.Cells.Item(7, 6).Value = 1.029629305299127
$ws.Cells.Item(7, 9).Value = 1.033957123087836
$ws.Cells.Item(7, 10).Value = 1.034591592358848
$ws.Cells.Item(7, 11).Value = 1.035551665748084
$ws.Cells.Item(7, 12).Value = 0.9968970624462087
$ws.Cells.Item(7, 13).Value = 1.032128792503456
$ws.Cells.Item(7, 14).Value = 1.015520491271521

$ws.Cells.Item(8, 2).Value = 1.02
$ws.Cells.Item(8, 3).Value = 1.028748760860455
$ws.Cells.Item(8, 4).Value = 1.032126462579505
$ws.Cells.Item(8, 5).Value = 0.9929600610674301
$ws.Cells.Item(8, 6).Value = 1.027513144531192
$ws.Cells.Item(8, 9).Value = 1.033658807057193
$ws.Cells.Item(8, 10).Value = 1.033779800662919
$ws.Cells.Item(8, 11).Value = 1.03486884846677
$ws.Cells.Item(8, 12).Value = 0.995817528259106
$ws.Cells.Item(8, 13).Value = 1.030268590407175
$ws.Cells.Item(8, 14).Value = 1.015248559645626

$ws.Cells.Item(9, 2).Value = 1.02
$ws.Cells.Item(9, 3).Value = 1.026481682416131
$ws.Cells.Item(9, 4).Value = 1.030476931819875
$ws.Cells.Item(9, 5).Value = 0.9906006454969559
$ws.Cells.Item(9, 6).Value = 1.02377985510899
$ws.Cells.Item(9, 9).Value = 1.033122539585371
$ws.Cells.Item(9, 10).Value = 1.032341240155077
$ws.Cells.Item(9, 11).Value = 1.033657632449482
$ws.Cells.Item(9, 12).Value = 0.9939188001724441
$ws.Cells.Item(9, 13).Value = 1.026982818756607
$ws.Cells.Item(9, 14).Value = 1.01476620029354

$ws.Cells.Item(10, 2).Value = 1.02
$ws.Cells.Item(10, 3).Value = 1.024967787736831
$ws.Cells.Item(10, 4).Value = 1.029375307484748
$ws.Cells.Item(10, 5).Value = 0.989033133672735
$ws.Cells.Item(10, 6).Value = 1.021287658611161
$ws.Cells.Item(10, 9).Value = 1.032757886624472
$ws.Cells.Item(10, 10).Value = 1.031376716880724
$ws.Cells.Item(10, 11).Value = 1.032844743491522
$ws.Cells.Item(10, 12).Value = 0.9926553831429383
$ws.Cells.Item(10, 13).Value = 1.02478665956159
$ws.Cells.Item(10, 14).Value = 1.014442475709074

$ws.Cells.Item(11, 2).Value = 1.02
$ws.Cells.Item(11, 3).Value = 1.024311642971531
$ws.Cells.Item(11, 4).Value = 1.028897828963703
$ws.Cells.Item(11, 5).Value = 0.988355674866747
$ws.Cells.Item(11, 6).Value = 1.020207625516931
$ws.Cells.Item(11, 9).Value = 1.032598289855339
$ws.Cells.Item(11, 10).Value = 1.030957755172423
$ws.Cells.Item(11, 11).Value = 1.032491462691158
$ws.Cells.Item(11, 12).Value = 0.9921088820399291
$ws.Cells.Item(11, 13).Value = 1.023834278628502
$ws.Cells.Item(11, 14).Value = 1.014301785872651

$ws.Cells.Item(12, 2).Value = 1.02
$ws.Cells.Item(12, 3).Value = 1.02406782714936
$ws.Cells.Item(12, 4).Value = 1.028720401470667
$ws.Cells.Item(12, 5).Value = 0.9881042295826724
$ws.Cells.Item(12, 6).Value = 1.019806310510694
$ws.Cells.Item(12, 9).Value = 1.032538752735177
$ws.Cells.Item(12, 10).Value = 1.03080193518326
$ws.Cells.Item(12, 11).Value = 1.032360043263364
$ws.Cells.Item(12, 12).Value = 0.9919059725120875
$ws.Cells.Item(12, 13).Value = 1.023480300136319
$ws.Cells.Item(12, 14).Value = 1.014249449728647

$ws.Cells.Item(13, 2).Value = 1.02
$ws.Cells.Item(13, 3).Value = 1.024120130769373
$ws.Cells.Item(13, 4).Value = 1.028758463476774
$ws.Cells.Item(13, 5).Value = 0.9881581567098651
$ws.Cells.Item(13, 6).Value = 1.019892400526162
$ws.Cells.Item(13, 9).Value = 1.032551535230728
$ws.Cells.Item(13, 10).Value = 1.030835368121308
$ws.Cells.Item(13, 11).Value = 1.032388242022054
$ws.Cells.Item(13, 12).Value = 0.9919494934313052
$ws.Cells.Item(13, 13).Value = 1.023556239875292
$ws.Cells.Item(13, 14).Value = 1.014260679529133

$ws.Cells.Item(14, 2).Value = 1.02
$ws.Cells.Item(14, 3).Value = 1.024291491015129
$ws.Cells.Item(14, 4).Value = 1.028883164196686
$ws.Cells.Item(14, 5).Value = 0.9883348863814464
$ws.Cells.Item(14, 6).Value = 1.020174455649751
$ws.Cells.Item(14, 9).Value = 1.032593373718596
$ws.Cells.Item(14, 10).Value = 1.030944879108727
$ws.Cells.Item(14, 11).Value = 1.032480603506674
$ws.Cells.Item(14, 12).Value = 0.9920921077337197
$ws.Cells.Item(14, 13).Value = 1.023805023221618
$ws.Cells.Item(14, 14).Value = 1.014297461336999

$ws.Cells.Item(15, 2).Value = 1.02
$ws.Cells.Item(15, 3).Value = 1.024397059157516
$ws.Cells.Item(15, 4).Value = 1.028959987043477
$ws.Cells.Item(15, 5).Value = 0.9884438009545853
$ws.Cells.Item(15, 6).Value = 1.020348219978563
$ws.Cells.Item(15, 9).Value = 1.032619117886511
$ws.Cells.Item(15, 10).Value = 1.031012326044804
$ws.Cells.Item(15, 11).Value = 1.032537484576829
$ws.Cells.Item(15, 12).Value = 0.9921799884222134
$ws.Cells.Item(15, 13).Value = 1.023958277249037
$ws.Cells.Item(15, 14).Value = 1.014320113520844

$ws.Cells.Item(16, 2).Value = 1.02
$ws.Cells.Item(16, 3).Value = 1.025011320744527
$ws.Cells.Item(16, 4).Value = 1.0294069862684
$ws.Cells.Item(16, 5).Value = 0.9890781214508737
$ws.Cells.Item(16, 6).Value = 1.021359317293613
$ws.Cells.Item(16, 9).Value = 1.032768442691097
$ws.Cells.Item(16, 10).Value = 1.03140449415081
$ws.Cells.Item(16, 11).Value = 1.032868162239355
$ws.Cells.Item(16, 12).Value = 0.9926916645766087
$ws.Cells.Item(16, 13).Value = 1.024849835229295
$ws.Cells.Item(16, 14).Value = 1.014451801954441

$ws.Cells.Item(17, 2).Value = 1.02
$ws.Cells.Item(17, 3).Value = 1.025396464171664
$ws.Cells.Item(17, 4).Value = 1.029687251512605
$ws.Cells.Item(17, 5).Value = 0.989476357848556
$ws.Cells.Item(17, 6).Value = 1.0219933064717
$ws.Cells.Item(17, 9).Value = 1.032861654859109
$ws.Cells.Item(17, 10).Value = 1.031650137515711
$ws.Cells.Item(17, 11).Value = 1.033075240515559
$ws.Cells.Item(17, 12).Value = 0.9930127773699352
$ws.Cells.Item(17, 13).Value = 1.025408698318744
$ws.Cells.Item(17, 14).Value = 1.014534268545484

$ws.Cells.Item(18, 2).Value = 1.02
$ws.Cells.Item(18, 3).Value = 1.02562105225589
$ws.Cells.Item(18, 4).Value = 1.029850680369281
$ws.Cells.Item(18, 5).Value = 0.9897087662937556
$ws.Cells.Item(18, 6).Value = 1.022363015768254
$ws.Cells.Item(18, 9).Value = 1.032915859978669
$ws.Cells.Item(18, 10).Value = 1.031793290162017
$ws.Cells.Item(18, 11).Value = 1.033195900917472
$ws.Cells.Item(18, 12).Value = 0.9932001317071769
$ws.Cells.Item(18, 13).Value = 1.02573453625107
$ws.Cells.Item(18, 14).Value = 1.014582320248063

$ws.Cells.Item(19, 2).Value = 1.02
$ws.Cells.Item(19, 3).Value = 1.025697620932742
$ws.Cells.Item(19, 4).Value = 1.029906397745451
$ws.Cells.Item(19, 5).Value = 0.9897880325774034
$ws.Cells.Item(19, 6).Value = 1.022489062765149
$ws.Cells.Item(19, 9).Value = 1.032934314725833
$ws.Cells.Item(19, 10).Value = 1.031842080012376
$ws.Cells.Item(19, 11).Value = 1.033237021840927
$ws.Cells.Item(19, 12).Value = 0.9932640239640975
$ws.Cells.Item(19, 13).Value = 1.025845615464869
$ws.Cells.Item(19, 14).Value = 1.014598696222179

$ws.Cells.Item(20, 2).Value = 1.02
$ws.Cells.Item(20, 3).Value = 1.025355148098368
$ws.Cells.Item(20, 4).Value = 1.029657186373927
$ws.Cells.Item(20, 5).Value = 0.9894336180360679
$ws.Cells.Item(20, 6).Value = 1.021925294371766
$ws.Cells.Item(20, 9).Value = 1.032851671037853
$ws.Cells.Item(20, 10).Value = 1.031623795448693
$ws.Cells.Item(20, 11).Value = 1.03305303589548
$ws.Cells.Item(20, 12).Value = 0.9929783193494215
$ws.Cells.Item(20, 13).Value = 1.025348751865707
$ws.Cells.Item(20, 14).Value = 1.014525425799937

$ws.Cells.Item(21, 2).Value = 1.02
$ws.Cells.Item(21, 3).Value = 1.024241032313255
$ws.Cells.Item(21, 4).Value = 1.028846444898088
$ws.Cells.Item(21, 5).Value = 0.9882828385668249
$ws.Cells.Item(21, 6).Value = 1.020091401448182
$ws.Cells.Item(21, 9).Value = 1.032581060390759
$ws.Cells.Item(21, 10).Value = 1.030912636353706
$ws.Cells.Item(21, 11).Value = 1.032453410746071
$ws.Cells.Item(21, 12).Value = 0.9920501090198102
$ws.Cells.Item(21, 13).Value = 1.023731768906641
$ws.Cells.Item(21, 14).Value = 1.014286632158041

$ws.Cells.Item(22, 2).Value = 1.02
$ws.Cells.Item(22, 3).Value = 1.023539996319677
$ws.Cells.Item(22, 4).Value = 1.028336290230541
$ws.Cells.Item(22, 5).Value = 0.9875604150241495
$ws.Cells.Item(22, 6).Value = 1.018937530954149
$ws.Cells.Item(22, 9).Value = 1.032409436693155
$ws.Cells.Item(22, 10).Value = 1.03046435073721
$ws.Cells.Item(22, 11).Value = 1.032075272832816
$ws.Cells.Item(22, 12).Value = 0.9914670000341481
$ws.Cells.Item(22, 13).Value = 1.022713820709474
$ws.Cells.Item(22, 14).Value = 1.01413604356069

$ws.Cells.Item(23, 2).Value = 1.02
$ws.Cells.Item(23, 3).Value = 1.023911681147356
$ws.Cells.Item(23, 4).Value = 1.028606771812583
$ws.Cells.Item(23, 5).Value = 0.9879432794643023
$ws.Cells.Item(23, 6).Value = 1.019549300846208
$ws.Cells.Item(23, 9).Value = 1.032500558096374
$ws.Cells.Item(23, 10).Value = 1.030702104968981
$ws.Cells.Item(23, 11).Value = 1.032275838229682
$ws.Cells.Item(23, 12).Value = 0.991776070289318
$ws.Cells.Item(23, 13).Value = 1.023253578648305
$ws.Cells.Item(23, 14).Value = 1.014215916142429

$ws.Cells.Item(24, 2).Value = 1.02
$ws.Cells.Item(24, 3).Value = 1.025373817237039
$ws.Cells.Item(24, 4).Value = 1.029670771655148
$ws.Cells.Item(24, 5).Value = 0.9894529299347244
$ws.Cells.Item(24, 6).Value = 1.021956026375627
$ws.Cells.Item(24, 9).Value = 1.032856182803388
$ws.Cells.Item(24, 10).Value = 1.031635698686872
$ws.Cells.Item(24, 11).Value = 1.033063069593025
$ws.Cells.Item(24, 12).Value = 0.9929938892766442
$ws.Cells.Item(24, 13).Value = 1.025375839510528
$ws.Cells.Item(24, 14).Value = 1.014529421609368

$ws.Cells.Item(25, 2).Value = 1.02
$ws.Cells.Item(25, 3).Value = 1.02706821359582
$ws.Cells.Item(25, 4).Value = 1.030903715604053
$ws.Cells.Item(25, 5).Value = 0.9912096547607049
$ws.Cells.Item(25, 6).Value = 1.024745557416192
$ws.Cells.Item(25, 9).Value = 1.03326243503634
$ws.Cells.Item(25, 10).Value = 1.032714105188843
$ws.Cells.Item(25, 11).Value = 1.033971712443208
$ws.Cells.Item(25, 12).Value = 0.9944092447426414
$ws.Cells.Item(25, 13).Value = 1.027833237803553
$ws.Cells.Item(25, 14).Value = 1.014891280035216
